# edit.ps1 - PowerPoint COM-interop (PowerShell) script
#
# Reproduces the two authored edits from the commit:
#   1. Slide 4 ("What is the output?") - content placeholder's second
#      paragraph: the leading word "Write" is retyped as "SPEF", turning
#      "Write File (Aya)" into "SPEF File (Aya)".
#   2. A new "Title and Content" slide titled "Progress Update " is
#      inserted right before the final "Test Cases Output" slide (i.e.
#      at position 5 of what becomes a 6-slide deck), with an empty
#      content placeholder.

$p = $ppt.ActivePresentation

# --- 1. Slide 4: "Write File (Aya)" -> "SPEF File (Aya)" -----------------
# Slide 4 ("What is the output?") uses the "Title and Content" layout:
# Placeholders.Item(1) is the title, Item(2) is the content body - find
# the content body defensively by name in case ordering ever shifts.
$slide4 = $p.Slides.Item(4)
$body = $null
for ($i = 1; $i -le $slide4.Shapes.Placeholders.Count; $i++) {
    $ph = $slide4.Shapes.Placeholders.Item($i)
    if ($ph.Name -like "Content Placeholder*") {
        $body = $ph.TextFrame.TextRange
    }
}
if ($null -eq $body) {
    $body = $slide4.Shapes.Placeholders.Item(2).TextFrame.TextRange
}

$secondPara = $body.Paragraphs(2)

# Replace the leading "Write " (including the trailing space) with
# "SPEF " so the rest of the line ("File (Aya)") is left untouched -
# matching the author's edit of retyping just the first word.
$leadWord = $secondPara.Characters(1, 6)
$leadWord.Text = "SPEF "

# --- 2. Insert new "Progress Update" slide before the last slide --------
# The deck currently ends with "Test Cases Output" (last slide); the new
# slide is inserted right before it and pushes it down by one position.
$insertAt = $p.Slides.Count
$newSlide = $p.Slides.Add($insertAt, 2)  # 2 = ppLayoutText (Title and Content)

$title = $null
for ($i = 1; $i -le $newSlide.Shapes.Placeholders.Count; $i++) {
    $ph = $newSlide.Shapes.Placeholders.Item($i)
    if ($ph.Name -like "Title*") {
        $title = $ph
    }
}
if ($null -eq $title) {
    $title = $newSlide.Shapes.Placeholders.Item(1)
}
$title.TextFrame.TextRange.Text = "Progress Update "
